$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1418.1428
$ws.Range("I19").Value = 4000
$ws.Range("K19").Value = 4000
$ws.Range("M19").Value = -3825

$ws.Range("H70").Value = 1475
$ws.Range("I70").Value = 900
$ws.Range("J70").Value = 1666.6666
$ws.Range("K70").Value = 2700
$ws.Range("L70").Value = 4999.9998
$ws.Range("M70").Value = -2430
$ws.Range("N70").Value = -5539.9998

$ws.Range("H73").Value = 1475
$ws.Range("I73").Value = 900
$ws.Range("J73").Value = 1666.6666
$ws.Range("K73").Value = 2700
$ws.Range("L73").Value = 4999.9998
$ws.Range("M73").Value = -1764
$ws.Range("N73").Value = -6871.9998

$ws.Range("H74").Value = 3600
$ws.Range("I74").Value = 3600
$ws.Range("K74").Value = 3600
$ws.Range("M74").Value = -2664

$ws.Range("H77").Value = 3600
$ws.Range("I77").Value = 3600
$ws.Range("K77").Value = 18000
$ws.Range("M77").Value = -13320

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2410.75
$ws.Range("I2").Value = 2473.6667
$ws.Range("J2").Value = 2222
$ws.Range("K2").Value = 2473.6667
$ws.Range("L2").Value = 2222
$ws.Range("M2").Value = -2360.6667
$ws.Range("N2").Value = -2448

$ws.Range("H44").Value = 80000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H116").Value = 2410.75
$ws.Range("I116").Value = 2473.6667
$ws.Range("J116").Value = 2222
$ws.Range("K116").Value = 2473.6667
$ws.Range("L116").Value = 2222
$ws.Range("M116").Value = -179.6667000000002
$ws.Range("N116").Value = -6810

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2410.75
$ws.Range("I3").Value = 2473.6667
$ws.Range("J3").Value = 2222
$ws.Range("K3").Value = 2473.6667
$ws.Range("L3").Value = 2222
$ws.Range("M3").Value = -2359.6667
$ws.Range("N3").Value = -2450

$ws.Range("H29").Value = 15596.2
$ws.Range("J29").Value = 9999.5
$ws.Range("L29").Value = 9999.5
$ws.Range("N29").Value = -10577.5

$ws.Range("H35").Value = 60000
$ws.Range("J35").Value = 60000
$ws.Range("L35").Value = 60000
$ws.Range("N35").Value = -60620

$ws.Range("H86").Value = 2511.2
$ws.Range("I86").Value = 2826.5
$ws.Range("J86").Value = 1250
$ws.Range("K86").Value = 2826.5
$ws.Range("L86").Value = 1250
$ws.Range("M86").Value = -1703.5
$ws.Range("N86").Value = -3496

$ws.Range("H89").Value = 2511.2
$ws.Range("I89").Value = 2826.5
$ws.Range("J89").Value = 1250
$ws.Range("K89").Value = 14132.5
$ws.Range("L89").Value = 6250
$ws.Range("M89").Value = -8516.5
$ws.Range("N89").Value = -17482

$ws.Range("H94").Value = 1122.8
$ws.Range("I94").Value = 1122.8
$ws.Range("K94").Value = 1122.8
$ws.Range("M94").Value = -671.8

$ws.Range("H99").Value = 2485
$ws.Range("I99").Value = 2485
$ws.Range("K99").Value = 2485
$ws.Range("M99").Value = -987

$ws.Range("H106").Value = 13685.875
$ws.Range("J106").Value = 13685.875
$ws.Range("L106").Value = 13685.875
$ws.Range("N106").Value = -16209.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3724.7
$ws.Range("I132").Value = 2655.875
$ws.Range("K132").Value = 7967.625
$ws.Range("M132").Value = -5437.625

$ws.Range("H134").Value = 2843.2
$ws.Range("I134").Value = 2843.2
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8529.599999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5994.599999999999
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2537.25
$ws.Range("I2").Value = 42.57143
$ws.Range("K2").Value = 255.42858
$ws.Range("M2").Value = -142.42858

$ws.Range("H4").Value = 1032.3334
$ws.Range("I4").Value = 261.25
$ws.Range("K4").Value = 783.75
$ws.Range("M4").Value = -671.75

$ws.Range("H17").Value = 132.77777
$ws.Range("I17").Value = 116.42857
$ws.Range("J17").Value = 190
$ws.Range("K17").Value = 349.28571
$ws.Range("L17").Value = 570
$ws.Range("M17").Value = -180.28571
$ws.Range("N17").Value = -908

$ws.Range("H25").Value = 395.6
$ws.Range("I25").Value = 366
$ws.Range("J25").Value = 440
$ws.Range("K25").Value = 1098
$ws.Range("L25").Value = 1320
$ws.Range("M25").Value = -929
$ws.Range("N25").Value = -1658

$ws.Range("H30").Value = 395.6
$ws.Range("I30").Value = 366
$ws.Range("J30").Value = 440
$ws.Range("K30").Value = 1098
$ws.Range("L30").Value = 1320
$ws.Range("M30").Value = -996
$ws.Range("N30").Value = -1524

$ws.Range("H34").Value = 1050.25
$ws.Range("J34").Value = 1171.7142
$ws.Range("L34").Value = 3515.1426
$ws.Range("N34").Value = -3683.1426

$ws.Range("H39").Value = 500
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 6000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H107").Value = 2388.25
$ws.Range("I107").Value = 1517.6666
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 4552.9998
$ws.Range("L107").Value = 15000
$ws.Range("M107").Value = -2632.9998
$ws.Range("N107").Value = -18840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2216.5
$ws.Range("I70").Value = 2216.5
$ws.Range("K70").Value = 2216.5
$ws.Range("M70").Value = -1946.5

$ws.Range("H73").Value = 2216.5
$ws.Range("I73").Value = 2216.5
$ws.Range("K73").Value = 2216.5
$ws.Range("M73").Value = -1280.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 3863
$ws.Range("I13").Value = 148.66667
$ws.Range("K13").Value = 148.66667
$ws.Range("M13").Value = -8.666670000000011

$ws.Range("H69").Value = 17802.334
$ws.Range("J69").Value = 17802.334
$ws.Range("L69").Value = 17802.334
$ws.Range("N69").Value = -19300.334

$ws.Range("H72").Value = 17802.334
$ws.Range("J72").Value = 17802.334
$ws.Range("L72").Value = 53407.00199999999
$ws.Range("N72").Value = -60895.00199999999

$ws.Range("H101").Value = 22180
$ws.Range("J101").Value = 22475
$ws.Range("L101").Value = 22475
$ws.Range("N101").Value = -28965

$ws.Range("H126").Value = 1409.5
$ws.Range("I126").Value = 1409.5
$ws.Range("K126").Value = 4228.5
$ws.Range("M126").Value = -1758.5
